# Renames the embedded logo pictures -- i.e. sets the Word "Name" of each
# InlineShape, which is backed by the <wp:docPr name="..."/> /
# <pic:cNvPr name="..."/> attributes on the drawing -- for:
#   * the Pearson Edexcel logo in the default (primary) footer
#   * the Pearson Edexcel logo in the first-page footer
#   * the BTEC logo in the first-page header
#
#   footer (default)    : PearsonLogo       image1.png -> image2.png
#   footer (first page) : PearsonLogo       image1.png -> image2.png
#   header (first page) : BTec_Logo-Orange  image2.jpg -> image1.jpg

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footer (default / wdHeaderFooterPrimary) ------------------------------
$ftrPrimary = $sec.Footers.Item(1)
if ($ftrPrimary.Exists -and $ftrPrimary.Range.InlineShapes.Count -ge 1) {
    $logo = $ftrPrimary.Range.InlineShapes.Item(1)
    $asShape = $logo.ConvertToShape()
    $asShape.Name = "image2.png"
    $asShape.ConvertToInlineShape() | Out-Null
}

# --- Footer (first page / wdHeaderFooterFirstPage) -------------------------
$ftrFirst = $sec.Footers.Item(2)
if ($ftrFirst.Exists -and $ftrFirst.Range.InlineShapes.Count -ge 1) {
    $logo = $ftrFirst.Range.InlineShapes.Item(1)
    $asShape = $logo.ConvertToShape()
    $asShape.Name = "image2.png"
    $asShape.ConvertToInlineShape() | Out-Null
}

# --- Header (first page / wdHeaderFooterFirstPage) --------------------------
$hdrFirst = $sec.Headers.Item(2)
if ($hdrFirst.Exists -and $hdrFirst.Range.InlineShapes.Count -ge 1) {
    $logo = $hdrFirst.Range.InlineShapes.Item(1)
    $logo.Name = "image1.jpg"
}
